$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.907.94"
$ws.Range("E2").Value = "  +2.46%  "

$ws.Range("D3").Value = "3.578.03"
$ws.Range("E3").Value = "  +2.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "624.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.82%  "

$ws.Range("D7").Value = "3.578.48"
$ws.Range("E7").Value = "  +2.57%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("E10").Value = "  +7.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.05%  "

$ws.Range("E13").Value = "  +4.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.82%  "

$ws.Range("D15").Value = "4.181.78"
$ws.Range("E15").Value = "  +2.46%  "

$ws.Range("D16").Value = "69.096.26"
$ws.Range("E16").Value = "  +2.93%  "

$ws.Range("D17").Value = "3.574.65"
$ws.Range("E17").Value = "  +2.70%  "

$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.61%  "

$ws.Range("E21").Value = "  +12.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "462.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.648"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("E25").Value = "  +4.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.85%  "

$ws.Range("D27").Value = "3.717.69"
$ws.Range("E27").Value = "  +2.40%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.63%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.00%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.170"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.27%  "

$ws.Range("D37").Value = "3.569.44"
$ws.Range("E37").Value = "  +2.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.10%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "179.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0926"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +16.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.905"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.25%  "

$ws.Range("E49").Value = "  +6.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.16%  "

$ws.Range("E51").Value = "  +9.36%  "
